$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts existing rows 23-46 down to 24-47.
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the new data record.
$ws.Range("A23").Value = 8
$ws.Range("B23").Value = "Terminal La Palmera de La Serena"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 45040
$ws.Range("D23").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100101
$ws.Range("H23").Value = "Berries"
$ws.Range("I23").Value = 100101001
$ws.Range("J23").Value = "Arándano (blue)"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 11000
$ws.Range("P23").Value = 10500
$ws.Range("Q23").Value = "$/bandeja 2 kilos"
$ws.Range("R23").Value = "Provincia de Curicó"
$ws.Range("S23").Value = 5250
$ws.Range("T23").Value = 2
